# Apply cell-level updates for existing rows (2-35) per the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44413
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 7000
$ws.Range("P2").Value = 583

$ws.Range("D3").Value = 44432

$ws.Range("D4").Value = 44259
$ws.Range("J4").Value = 30

$ws.Range("D5").Value = 44186
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = 5000
$ws.Range("P5").Value = 417

$ws.Range("D6").Value = 44326
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 6000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 6000
$ws.Range("P6").Value = 500

$ws.Range("D7").Value = 44428
$ws.Range("J7").Value = 10
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 7000
$ws.Range("P7").Value = 583

$ws.Range("D8").Value = 44454

$ws.Range("D9").Value = 44424
$ws.Range("J9").Value = 30
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 7000
$ws.Range("P9").Value = 583

$ws.Range("D10").Value = 44179
$ws.Range("K10").Value = 6000
$ws.Range("L10").Value = 6000
$ws.Range("M10").Value = 6000
$ws.Range("P10").Value = 500

$ws.Range("D11").Value = 44329
$ws.Range("J11").Value = 40
$ws.Range("M11").Value = 5500
$ws.Range("P11").Value = 458

$ws.Range("D12").Value = 44302
$ws.Range("J12").Value = 20

$ws.Range("D13").Value = 44435
$ws.Range("J13").Value = 30
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 7000
$ws.Range("P13").Value = 583

$ws.Range("D14").Value = 44195
$ws.Range("J14").Value = 55
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = 5000
$ws.Range("P14").Value = 417

$ws.Range("D15").Value = 44449
$ws.Range("J15").Value = 65
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 7000
$ws.Range("P15").Value = 583

$ws.Range("D17").Value = 44327
$ws.Range("J17").Value = 30
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = 6000
$ws.Range("P17").Value = 500

$ws.Range("D18").Value = 44452
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 7000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 7000
$ws.Range("P18").Value = 583

$ws.Range("D19").Value = 44441

$ws.Range("D20").Value = 44442
$ws.Range("K20").Value = 6000
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = 6000
$ws.Range("P20").Value = 500

$ws.Range("D21").Value = 44453
$ws.Range("J21").Value = 20

$ws.Range("D22").Value = 44448
$ws.Range("J22").Value = 50

$ws.Range("D23").Value = 44165
$ws.Range("J23").Value = 130
$ws.Range("K23").Value = 5000
$ws.Range("M23").Value = 5615
$ws.Range("P23").Value = 468

$ws.Range("D24").Value = 44427
$ws.Range("J24").Value = 20
$ws.Range("K24").Value = 7000
$ws.Range("L24").Value = 7000
$ws.Range("M24").Value = 7000
$ws.Range("P24").Value = 583

$ws.Range("D25").Value = 44196
$ws.Range("J25").Value = 20
$ws.Range("K25").Value = 5000
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = 5000
$ws.Range("P25").Value = 417

$ws.Range("D26").Value = 44301
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 6000
$ws.Range("L26").Value = 6000
$ws.Range("M26").Value = 6000
$ws.Range("O26").Value = "Provincia de Cautín"
$ws.Range("P26").Value = 500

$ws.Range("D27").Value = 44166
$ws.Range("J27").Value = 55

$ws.Range("D28").Value = 44369
$ws.Range("K28").Value = 4000
$ws.Range("L28").Value = 4000
$ws.Range("M28").Value = 4000
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 333

$ws.Range("D29").Value = 44438
$ws.Range("J29").Value = 30
$ws.Range("K29").Value = 6000
$ws.Range("M29").Value = 6000
$ws.Range("P29").Value = 500

$ws.Range("D30").Value = 44372
$ws.Range("J30").Value = 40
$ws.Range("K30").Value = 6000
$ws.Range("L30").Value = 6000
$ws.Range("M30").Value = 6000
$ws.Range("P30").Value = 500

$ws.Range("D31").Value = 44203
$ws.Range("J31").Value = 40
$ws.Range("K31").Value = 5000
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = 5000
$ws.Range("P31").Value = 417

$ws.Range("D32").Value = 44162
$ws.Range("J32").Value = 50
$ws.Range("K32").Value = 5000
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = 5000
$ws.Range("P32").Value = 417

$ws.Range("D33").Value = 44410
$ws.Range("J33").Value = 40
$ws.Range("K33").Value = 7000
$ws.Range("L33").Value = 7000
$ws.Range("M33").Value = 7000
$ws.Range("P33").Value = 583

$ws.Range("D34").Value = 44411
$ws.Range("J34").Value = 20

$ws.Range("D35").Value = 44211
$ws.Range("J35").Value = 65
$ws.Range("K35").Value = 5000
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = 5000
$ws.Range("P35").Value = 417

# Add new row 36 (new data record)
$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = 44425
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = 300000001
$ws.Range("G36").Value = "Rabanito"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 20
$ws.Range("K36").Value = 7000
$ws.Range("L36").Value = 7000
$ws.Range("M36").Value = 7000
$ws.Range("N36").Value = "$/docena de paquetes"
$ws.Range("O36").Value = "Provincia de Cautín"
$ws.Range("P36").Value = 583
$ws.Range("Q36").Value = 12
$ws.Range("R36").Value = "Hortaliza"

# Apply date number format to the new date cell to match column D style
$ws.Range("D36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
